# Apply updated odds/statistics values to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7
$ws.Range("G7").Value = 2.25
$ws.Range("I7").Value = 2.9
$ws.Range("J7").Value = 2.88
$ws.Range("K7").Value = 2.38
$ws.Range("L7").Value = 3.25
$ws.Range("U7").Value = 1.5
$ws.Range("V7").Value = 2.5
$ws.Range("X7").Value = 13
$ws.Range("Y7").Value = 9.5
$ws.Range("AA7").Value = 17
$ws.Range("AE7").Value = 11
$ws.Range("AJ7").Value = 29
$ws.Range("AL7").Value = 23
$ws.Range("AM7").Value = 101
$ws.Range("AO7").Value = 12
$ws.Range("AW7").Value = 5
$ws.Range("AZ7").Value = 41
$ws.Range("BB7").Value = 101

# Row 8
$ws.Range("AE8").Value = 12
$ws.Range("AV8").Value = 41

# Row 14
$ws.Range("S14").Value = 1.3
$ws.Range("T14").Value = 3.4
$ws.Range("W14").Value = 9.5
$ws.Range("AL14").Value = 34
$ws.Range("AT14").Value = 3.4
$ws.Range("BB14").Value = 151

# Row 17
$ws.Range("G17").Value = 3.25
$ws.Range("I17").Value = 2.25
$ws.Range("J17").Value = 4
$ws.Range("L17").Value = 3
$ws.Range("AQ17").Value = 67
$ws.Range("AR17").Value = 101
